$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 460, shifting existing rows
# (old 460-472) down to (462-474) and preserving their contents/styles.
$ws.Rows("460:461").Insert()

# Populate the newly inserted row 460 with the new weekly record.
$ws.Range("A460").Value = 11
$ws.Range("B460").Value = "Vega Monumental Concepción"
$ws.Range("C460").Value = "Bíobío"
$ws.Range("D460").Value = 44706
$ws.Range("E460").Value = 8
$ws.Range("F460").Value = 100112004
$ws.Range("G460").Value = "Cebolla"
$ws.Range("H460").Value = "Sin especificar"
$ws.Range("I460").Value = "1a (cosecha)"
$ws.Range("J460").Value = 300
$ws.Range("K460").Value = 6500
$ws.Range("L460").Value = 6500
$ws.Range("M460").Value = 6500
$ws.Range("N460").Value = "`$/malla 18 kilos"
$ws.Range("O460").Value = "Región de O'Higgins"
$ws.Range("P460").Value = 361
$ws.Range("Q460").Value = 18
$ws.Range("R460").Value = "Hortaliza"

# Populate the newly inserted row 461 with the new weekly record.
$ws.Range("A461").Value = 11
$ws.Range("B461").Value = "Vega Monumental Concepción"
$ws.Range("C461").Value = "Bíobío"
$ws.Range("D461").Value = 44706
$ws.Range("E461").Value = 8
$ws.Range("F461").Value = 100112004
$ws.Range("G461").Value = "Cebolla"
$ws.Range("H461").Value = "Sin especificar"
$ws.Range("I461").Value = "2a (cosecha)"
$ws.Range("J461").Value = 300
$ws.Range("K461").Value = 5500
$ws.Range("L461").Value = 5500
$ws.Range("M461").Value = 5500
$ws.Range("N461").Value = "`$/malla 18 kilos"
$ws.Range("O461").Value = "Región de O'Higgins"
$ws.Range("P461").Value = 306
$ws.Range("Q461").Value = 18
$ws.Range("R461").Value = "Hortaliza"
